$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$insertPoint = $titlePara.Range
$insertPoint.InsertParagraphAfter()

$newParaRange = $d.Paragraphs(2).Range
$newParaRange.Collapse(1)

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:r/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
    '<w:r><w:t>: Experience the stunning visuals and rewarding bonus features of Arctic Valor, an online slot game inspired by Norse mythology. Play for free and read our review here.</w:t></w:r>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

[void]$newParaRange.InsertXML($metaXml)

# ------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Arctic Valor..." paragraph that
#    used to sit just before the final italic paragraph at the very
#    end of the document.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($count - 1)
$followingPara = $d.Paragraphs($count)

if ($dupTitlePara.Range.Text.Contains("Play Arctic Valor Online Slot Game for Free | Review")) {
    $dupRange = $d.Range($dupTitlePara.Range.Start, $followingPara.Range.Start)
    $dupRange.Delete()
}

# ------------------------------------------------------------------
# 3) Replace the closing italic paragraph's text with the new image
#    generation prompt, keeping its italic formatting intact. Scope
#    the search to just that final paragraph so the identical phrase
#    inside the new Meta description paragraph (inserted in step 1)
#    is left untouched.
# ------------------------------------------------------------------
$oldText = "Experience the stunning visuals and rewarding bonus features of Arctic Valor, an online slot game inspired by Norse mythology. Play for free and read our review here."
$newText = "Prompt: Create a feature image for Arctic Valor that features a happy Maya warrior with glasses in a cartoon style. The main colors used should be blue and white to match the icy theme of the game. The warrior should be holding a shield with a precious gemstone at the center, and in the background, there should be swirling snow and icicles hanging from the top. The image should be action-packed and showcase the excitement of the game."

$finalParaRange = $d.Paragraphs($d.Paragraphs.Count).Range
[void]$finalParaRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
